{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst newText =\n  \"The two architectures for version management tools are Centralised systems and Distributed systems.\";\n\n// Find the \"A cool quote by Dijkstra:\" paragraph; the following paragraph\n// holds the quote itself. Together they collapse into a single paragraph\n// carrying the new sentence.\nlet introIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"A cool quote by Dijkstra:\") {\n    introIndex = i;\n    break;\n  }\n}\n\nif (introIndex !== -1 && introIndex + 1 < paragraphs.items.length) {\n  const introPara = paragraphs.items[introIndex];\n  const quotePara = paragraphs.items[introIndex + 1];\n  quotePara.delete();\n  introPara.insertText(newText, \"Replace\");\n} else {\n  // Fallback: locate the quote paragraph directly via search and replace its\n  // whole paragraph text, removing the intro paragraph if still present.\n  const results = body.search(\"Computer science is no more about computers\", {\n    matchCase: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    const para = results.items[0].paragraphs.getFirst();\n    para.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$newText = \"The two architectures for version management tools are Centralised systems and Distributed systems.\"\n\n# Locate the intro paragraph (\"A cool quote by Dijkstra:\"); the paragraph\n# immediately after it holds the Dijkstra quote. The two are merged into a\n# single paragraph carrying the new sentence.\n$introIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n  $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)\n  if ($t -eq \"A cool quote by Dijkstra:\") {\n    $introIndex = $i\n    break\n  }\n}\n\nif ($introIndex -ne -1 -and ($introIndex + 1) -le $d.Paragraphs.Count) {\n  $introPara = $d.Paragraphs.Item($introIndex)\n  $quotePara = $d.Paragraphs.Item($introIndex + 1)\n  $quotePara.Range.Delete()\n  $introPara.Range.Text = $newText\n} else {\n  # Fallback: find the quote text directly and replace the whole paragraph.\n  $rng = $d.Content\n  if ($rng.Find.Execute(\"Computer science is no more about computers\")) {\n    $rng.Expand(4) | Out-Null\n    $rng.Text = $newText\n  }\n}\n"}
